# The document contains a repeating set of "plumbing" paragraphs (styled
# NodeStart / NodeEnd / HeadEnd / etc.) used as separators around each
# resource. A bug produced extra blank "Node End" / "Head End" paragraphs
# that swallowed/duplicated header text around link resources. This script
# removes the superfluous separator paragraphs, restoring a single
# separator of each kind between sections.
#
# Paragraphs are identified by (1-indexed) position together with their
# style name and content as a safety check, then removed from the end of
# the document backwards so earlier indices stay valid as we go.

$d = $word.ActiveDocument

# (index, expected style, expected text) -- validated against the document
# before deleting anything.
$targets = @(
    @{ Index = 63; Style = "Node End"; Text = [string][char]0x00A0 },
    @{ Index = 55; Style = "Node End"; Text = [string][char]0x00A0 },
    @{ Index = 54; Style = "Node End"; Text = [string][char]0x00A0 },
    @{ Index = 52; Style = "Head End"; Text = [string][char]0x00A0 },
    @{ Index = 47; Style = "Node End"; Text = [string][char]0x00A0 },
    @{ Index = 46; Style = "Node End"; Text = [string][char]0x00A0 },
    @{ Index = 32; Style = "Head End"; Text = [string][char]0x00A0 },
    @{ Index = 27; Style = "Node End"; Text = [string][char]0x00A0 }
)

foreach ($t in $targets) {
    $p = $d.Paragraphs.Item($t.Index)
    $styleName = $p.Style.Name
    $text = $p.Range.Text.TrimEnd([char]0x0007, "`r", "`n")
    if ($styleName -eq $t.Style -and $text -eq $t.Text) {
        $p.Range.Delete()
    } else {
        Write-Host "Skipping index" $t.Index "- expected" $t.Style "got" $styleName "text[" $text "]"
    }
}
